# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed May 22 13:24:43 UTC 2024 with GitHub Actions"
#
# Every assigned literal is prefixed with a leading apostrophe so Excel
# stores it verbatim as text (matching the original inline-string cells)
# instead of auto-coercing number-looking values (e.g. "615.61", "0.530",
# "39.70") into a floating point number and losing formatting / zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.883.28"
$ws.Range("E2").Value = "'  -1.53%  "

$ws.Range("D3").Value = "'3.699.03"
$ws.Range("E3").Value = "'  -2.23%  "

$ws.Range("E4").Value = "'  -0.07%  "

$ws.Range("D5").Value = "'615.61"
$ws.Range("E5").Value = "'  +0.69%  "

$ws.Range("D6").Value = "'178.11"
$ws.Range("E6").Value = "'  -0.35%  "

$ws.Range("D7").Value = "'3.696.54"
$ws.Range("E7").Value = "'  -2.20%  "

$ws.Range("E8").Value = "'  -0.06%  "

$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "'  -2.16%  "

$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "'  -1.70%  "

$ws.Range("D11").Value = "'6.25"
$ws.Range("E11").Value = "'  -2.32%  "

$ws.Range("D12").Value = "'0.479"
$ws.Range("E12").Value = "'  -4.01%  "

$ws.Range("D13").Value = "'39.70"
$ws.Range("E13").Value = "'  -2.32%  "

$ws.Range("D14").Value = "'0.0000252"
$ws.Range("E14").Value = "'  -1.77%  "

$ws.Range("D15").Value = "'4.314.47"
$ws.Range("E15").Value = "'  -2.22%  "

$ws.Range("D16").Value = "'3.693.67"
$ws.Range("E16").Value = "'  -2.27%  "

$ws.Range("D17").Value = "'69.834.90"
$ws.Range("E17").Value = "'  -1.86%  "

$ws.Range("E18").Value = "'  -1.97%  "

$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "'  -0.26%  "

$ws.Range("D20").Value = "'16.33"
$ws.Range("E20").Value = "'  -2.33%  "

$ws.Range("D21").Value = "'500.34"
$ws.Range("E21").Value = "'  -4.33%  "

$ws.Range("D22").Value = "'9.14"
$ws.Range("E22").Value = "'  -3.28%  "

$ws.Range("D23").Value = "'0.711"
$ws.Range("E23").Value = "'  -4.35%  "

$ws.Range("D24").Value = "'2.55"
$ws.Range("E24").Value = "'  +2.86%  "

$ws.Range("D25").Value = "'86.10"
$ws.Range("E25").Value = "'  -2.60%  "

$ws.Range("D26").Value = "'11.45"
$ws.Range("E26").Value = "'  +3.90%  "

$ws.Range("D27").Value = "'12.93"
$ws.Range("E27").Value = "'  -4.31%  "

$ws.Range("D28").Value = "'0.0000128"
$ws.Range("E28").Value = "'  +5.72%  "

$ws.Range("E29").Value = "'  +0.20%  "

$ws.Range("E30").Value = "'  -3.11%  "

$ws.Range("D31").Value = "'2.89"
$ws.Range("E31").Value = "'  -0.80%  "

$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "'  -1.76%  "

$ws.Range("D33").Value = "'30.12"
$ws.Range("E33").Value = "'  -6.61%  "

$ws.Range("E34").Value = "'  -1.30%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "'  -0.19%  "

$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "'  -0.89%  "

$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = "'  -1.66%  "

$ws.Range("D38").Value = "'0.137"
$ws.Range("E38").Value = "'  +3.96%  "

$ws.Range("D39").Value = "'0.338"
$ws.Range("E39").Value = "'  -0.49%  "

$ws.Range("D40").Value = "'2.06"
$ws.Range("E40").Value = "'  -7.56%  "

$ws.Range("D41").Value = "'49.99"
$ws.Range("E41").Value = "'  -2.87%  "

$ws.Range("D42").Value = "'44.90"
$ws.Range("E42").Value = "'  +1.87%  "

$ws.Range("D43").Value = "'430.55"
$ws.Range("E43").Value = "'  +0.63%  "

$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "'  +3.35%  "

$ws.Range("D45").Value = "'8.53"
$ws.Range("E45").Value = "'  -3.21%  "

$ws.Range("D46").Value = "'2.945.35"
$ws.Range("E46").Value = "'  -6.61%  "

$ws.Range("D47").Value = "'0.0358"
$ws.Range("E47").Value = "'  -2.29%  "

$ws.Range("B48").Value = "'InjectiveProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'27.31"
$ws.Range("E48").Value = "'  -1.62%  "

$ws.Range("B49").Value = "'USDe"
$ws.Range("C49").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "'  +0.00%  "

$ws.Range("D50").Value = "'136.03"
$ws.Range("E50").Value = "'  -3.83%  "

$ws.Range("D51").Value = "'2.42"
$ws.Range("E51").Value = "'  -2.17%  "
